# Check Input Data.xlsx - update to split the combined BVTQaZ and VTQaZ
# transportation CSV entries (on the "Boolean" sheet) into their per-vehicle-type
# file names, matching the new eps-us InputData file layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boolean")

# --- Expand "trans/BVTQaZ/BVTQaZ.csv" (row 17) into 6 rows ---------------
$bvtqazRow = 0
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    if ($ws.Cells.Item($r, 1).Value2 -eq "trans/BVTQaZ/BVTQaZ.csv") {
        $bvtqazRow = $r
        break
    }
}

if ($bvtqazRow -gt 0) {
    # insert 5 extra rows right after the existing one, so we end up with 6 rows
    $ws.Rows.Item($bvtqazRow + 1).Resize(5).Insert()

    $bvtqazValues = @(
        "trans/BVTQaZ/BVTQaZ-LDVs.csv",
        "trans/BVTQaZ/BVTQaZ-HDVs.csv",
        "trans/BVTQaZ/BVTQaZ-aircraft.csv",
        "trans/BVTQaZ/BVTQaZ-rail.csv",
        "trans/BVTQaZ/BVTQaZ-ships.csv",
        "trans/BVTQaZ/BVTQaZ-motorbikes.csv"
    )

    for ($i = 0; $i -lt $bvtqazValues.Length; $i++) {
        $ws.Cells.Item($bvtqazRow + $i, 1).Value = $bvtqazValues[$i]
    }
}

# --- Expand "trans/VTQaZ/VTQaZ.csv" into 6 rows ---------------------------
$vtqazRow = 0
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    if ($ws.Cells.Item($r, 1).Value2 -eq "trans/VTQaZ/VTQaZ.csv") {
        $vtqazRow = $r
        break
    }
}

if ($vtqazRow -gt 0) {
    $ws.Rows.Item($vtqazRow + 1).Resize(5).Insert()

    $vtqazValues = @(
        "trans/VTQaZ/VTQaZ-LDVs.csv",
        "trans/VTQaZ/VTQaZ-HDVs.csv",
        "trans/VTQaZ/VTQaZ-aircraft.csv",
        "trans/VTQaZ/VTQaZ-rail.csv",
        "trans/VTQaZ/VTQaZ-ships.csv",
        "trans/VTQaZ/VTQaZ-motorbikes.csv"
    )

    for ($i = 0; $i -lt $vtqazValues.Length; $i++) {
        $ws.Cells.Item($vtqazRow + $i, 1).Value = $vtqazValues[$i]
    }
}

# --- Restore the handful of trailing blank (but styled) rows left over
#     below the last entry, matching the row formatting used for the list
#     above (rows 33:38 on the Boolean sheet) ---------------------------
$ws.Cells.Item(32, 1).Copy()
$blankRows = $ws.Range($ws.Cells.Item(33, 1), $ws.Cells.Item(38, 1))
$blankRows.PasteSpecial(-4122)  # xlPasteFormats
$blankRows.ClearContents()
$excel.CutCopyMode = 0

# --- View bookkeeping -----------------------------------------------------
# "Boolean" sheet: select the last entry (A32).
[void]$ws.Activate()
[void]$ws.Cells.Item(32, 1).Select()

# "Integer" sheet: select A13 (no longer the active tab).
$intSheet = $wb.Worksheets.Item("Integer")
[void]$intSheet.Activate()
[void]$intSheet.Range("A13").Select()

# "About" sheet becomes the active sheet/tab on reopen.
$about = $wb.Worksheets.Item("About")
[void]$about.Activate()
[void]$about.Range("A1").Select()
